$wb = $excel.ActiveWorkbook
$targetSheets = @("展览", "全部类型")
foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 2
    $ws.Range("C2").Value = '苏州·世纪幻想动漫游戏展'
    $ws.Range("F2").Value = 1974
    $ws.Range("G2").Value = 6000

    # Row 3
    $ws.Range("C3").Value = '昆山·“不是！你有病吧！”主题展（取消）'
    $ws.Range("G3").Value = 0

    # Row 4
    $ws.Range("C4").Value = '苏州.第二届THO 赤维极陵'
    $ws.Range("G4").Value = 5800

    # Row 5
    $ws.Range("C5").Value = '【会员购严选】苏州·二次元开放式年会- I COME ACG'
    $ws.Range("D5").Value = '金山南路288号木渎影视城F2 苏州广电国际会展中心'
    $ws.Range("E5").Value = '2024.02.03 10:00-02.03 20:00'
    $ws.Range("F5").Value = 10454
    $ws.Range("G5").Value = 2500
    $ws.Range("I5").Value = 'https://show.bilibili.com/platform/detail.html?id=80426'
    $ws.Range("J5").Value = '//i2.hdslb.com/bfs/openplatform/202401/IkyhIHPT1704352086775.jpeg'

    # Row 6
    $ws.Range("C6").Value = '苏州·TCD国潮动漫游戏嘉年华'
    $ws.Range("D6").Value = '苏州大道东688号 苏州国际博览中心'
    $ws.Range("E6").Value = '2024.02.03 09:30-02.04 17:00'
    $ws.Range("F6").Value = 9069
    $ws.Range("G6").Value = 6000
    $ws.Range("I6").Value = 'https://show.bilibili.com/platform/detail.html?id=80084'
    $ws.Range("J6").Value = '//i0.hdslb.com/bfs/openplatform/202401/aDe3s9MS1705479547745.jpeg'

    # Row 7
    $ws.Range("C7").Value = '苏州·TCD国潮动漫游戏嘉年华吴磊内场'
    $ws.Range("G7").Value = 0

    # Row 8
    $ws.Range("C8").Value = '常熟·CDW·动漫展02'
    $ws.Range("D8").Value = '常熟国际展览中心 国际展览中心'
    $ws.Range("E8").Value = '2024.02.14 09:00-02.15 17:30'
    $ws.Range("F8").Value = 663
    $ws.Range("G8").Value = 5500
    $ws.Range("H8").Value = $False
    $ws.Range("I8").Value = 'https://show.bilibili.com/platform/detail.html?id=80504'
    $ws.Range("J8").Value = '//i1.hdslb.com/bfs/openplatform/202401/VHHzVjad1704438989848.jpeg'

    # Row 9
    $ws.Range("C9").Value = '常熟·漫魂动漫游戏展01'
    $ws.Range("D9").Value = '虞山北路258号 星程酒店(长江路店)'
    $ws.Range("E9").Value = '2024.02.14 09:00-02.14 21:00'
    $ws.Range("F9").Value = 89
    $ws.Range("G9").Value = 5000
    $ws.Range("I9").Value = 'https://show.bilibili.com/platform/detail.html?id=80248'
    $ws.Range("J9").Value = '//i2.hdslb.com/bfs/openplatform/202312/oPrKUOby1703664065719.jpeg'

    # Row 10
    $ws.Range("D10").Value = '润元路润南巷172号,地铁二号线陆慕站东200米,近市旅游换乘中心北100米 斐利酒店'
    $ws.Range("E10").Value = '2024.02.14 10:00-02.14 16:00'
    $ws.Range("F10").Value = 21
    $ws.Range("G10").Value = 4900
    $ws.Range("H10").Value = $True
    $ws.Range("I10").Value = 'https://show.bilibili.com/platform/detail.html?id=80528'
    $ws.Range("J10").Value = '//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg'

    # Row 11
    $ws.Range("C11").Value = '太仓·龙狮新春动漫节4.0'
    $ws.Range("D11").Value = '滨河路126号 凯景世纪大酒店'
    $ws.Range("E11").Value = '2024.02.16 08:30-02.16 15:00'
    $ws.Range("F11").Value = 12
    $ws.Range("G11").Value = 4500
    $ws.Range("I11").Value = 'https://show.bilibili.com/platform/detail.html?id=81044'
    $ws.Range("J11").Value = '//i1.hdslb.com/bfs/openplatform/202401/AMDXVltp1705568031796.jpeg'

    # Row 12
    $ws.Range("C12").Value = '苏州·Good Jump ACG迎新特别篇X动漫品牌博览会'
    $ws.Range("D12").Value = '金山南路288号 广电国际会展中心'
    $ws.Range("E12").Value = '2024.02.16 10:00-02.17 17:00'
    $ws.Range("F12").Value = 9396
    $ws.Range("G12").Value = 6000
    $ws.Range("I12").Value = 'https://show.bilibili.com/platform/detail.html?id=79303'
    $ws.Range("J12").Value = '//i2.hdslb.com/bfs/openplatform/202312/C3P0Encm1701659824998.jpeg'

    # Row 13
    $ws.Range("C13").Value = '苏州·第五届次元鹿角动漫游戏展'
    $ws.Range("G13").Value = 6800

    # Row 14
    $ws.Range("C14").Value = '苏州·国风宠物-cosplay展'
    $ws.Range("G14").Value = 6500

    # Row 15
    $ws.Range("C15").Value = '苏州·绘时国乙1.0-秩序之外'
    $ws.Range("F15").Value = 56
    $ws.Range("G15").Value = 7800

    # Row 16
    $ws.Range("C16").Value = '苏州·梦幻岛 国乙主题文化展（日夜场） 梦幻岛之约3.0'
    $ws.Range("G16").Value = 4830

    # Row 17
    $ws.Range("C17").Value = '昆山·第十二届理想乡动漫游戏展'
    $ws.Range("D17").Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
    $ws.Range("E17").Value = '2024.05.01 10:00-05.03 17:00'
    $ws.Range("F17").Value = 10604
    $ws.Range("G17").Value = 10
    $ws.Range("I17").Value = 'https://show.bilibili.com/platform/detail.html?id=77196'
    $ws.Range("J17").Value = '//i2.hdslb.com/bfs/openplatform/202310/9xMTQMlg1696736126094.png'

    # Row 18
    $ws.Range("C18").Value = '苏州·第十七届 I COME ACG  动漫品牌博览会'
    $ws.Range("D18").Value = '金山南路288号 广电国际会展中心'
    $ws.Range("E18").Value = '2024.05.01 10:00-05.02 17:00'
    $ws.Range("F18").Value = 10670
    $ws.Range("G18").Value = 6500
    $ws.Range("I18").Value = 'https://show.bilibili.com/platform/detail.html?id=79789'
    $ws.Range("J18").Value = '//i2.hdslb.com/bfs/openplatform/202312/lau3mW031702535438289.jpeg'

    # Row 19
    $ws.Range("F19").Value = 3
    $ws.Range("G19").Value = 100
    $ws.Range("I19").Value = 'https://show.bilibili.com/platform/detail.html?id=81116'
    $ws.Range("J19").Value = '//i2.hdslb.com/bfs/openplatform/202401/EubrAneC1705648695005.jpeg'

    # Row 20
    $ws.Range("C20").Value = '昆山·第十二届理想乡动漫游戏展嘉宾啊川签售会'
    $ws.Range("G20").Value = 100
    $ws.Range("I20").Value = 'https://show.bilibili.com/platform/detail.html?id=81100'
    $ws.Range("J20").Value = '//i2.hdslb.com/bfs/openplatform/202401/F24i5GMX1705646667852.jpeg'

    # Row 21
    $ws.Range("C21").Value = '昆山·第十二届理想乡动漫游戏展嘉宾漠小然签售会'
    $ws.Range("G21").Value = 100

    # Row 22
    $ws.Range("C22").Value = '昆山·第十二届理想乡动漫游戏展嘉宾葫芦岛老八签售会'
    $ws.Range("F22").Value = 2
    $ws.Range("G22").Value = 100
    $ws.Range("I22").Value = 'https://show.bilibili.com/platform/detail.html?id=81118'
    $ws.Range("J22").Value = '//i2.hdslb.com/bfs/openplatform/202401/uHOCneLv1705648779163.jpeg'

    # Row 23
    $ws.Range("C23").Value = '昆山·第十二届理想乡动漫游戏展嘉宾沈辞签售会'
    $ws.Range("G23").Value = 100

    # Row 24
    $ws.Range("C24").Value = '昆山·第十二届理想乡动漫游戏展嘉宾矮乐多aliga签售会'
    $ws.Range("G24").Value = 100
}
